# Fix time unit for "complete" parameter: change from frame-count based
# ("最大补全帧数" / "maxCompleteFrames" / 20) to time based
# ("最大补全时间" / "maxCompleteTime(s)" / "2s").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update parameter name (C11) and symbol (D11) text.
$ws.Range("C11").Value = "最大补全时间"
$ws.Range("D11").Value = "maxCompleteTime(s)"

# Update the value cell (E11) from a numeric frame count to a textual
# time value ("2s" is not numeric, so Excel stores it as text).
$ws.Range("E11").Value = "2s"

# Reflect the new active cell selection recorded in the saved workbook.
$ws.Range("E12").Select()

$wb.Save()
